# DES 2024_evaluation.docx — apply the author's edits:
#  1. "By:  your name" -> "By:  Luke Hammond" (plain, non-italic run)
#  2. Replace the 3 italic placeholder lines under "Client details:" with
#     the real contact info (incl. two real hyperlinks) and drop one of
#     the two blank paragraphs that used to follow them.
#  3. Materialise the "Hyperlink" / "Unresolved Mention" character styles
#     that Word writes out once a real hyperlink is inserted.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "your name" -> "Luke Hammond" (drop the italic formatting)
# ---------------------------------------------------------------------
$byPara = $d.Paragraphs(6)
$find = $byPara.Range.Duplicate
$find.Find.ClearFormatting()
$find.Find.Text = "your name"
$find.Find.Execute() | Out-Null
$find.Delete()
$find.InsertAfter("Luke Hammond")

# ---------------------------------------------------------------------
# 2) Client-details block
# ---------------------------------------------------------------------
# Paragraphs 10-12 hold the three italic placeholder lines; delete their
# text plus the two trailing paragraph marks so we fall back to clean,
# unformatted paragraphs (paragraph 9, "Client details:", carries no
# special run/paragraph formatting).
$clientDetails = $d.Paragraphs(9)
$lastPlaceholder = $d.Paragraphs(12)
$d.Range($clientDetails.Range.End, $lastPlaceholder.Range.End).Delete()

# -- paragraph: "Kathryn Shatford - <email>" --
$p1 = $d.Paragraphs(10)
$d.Range($p1.Range.Start, $p1.Range.Start).InsertAfter("Kathryn Shatford - X ")
$p1 = $d.Paragraphs(10)
$xPos = $p1.Range.Start + 19
$emailSel = $d.Range($xPos, $xPos + 1)
$d.Hyperlinks.Add($emailSel, "mailto:kathryn.shatford@housingmatters.org.uk", $null, $null, "kathryn.shatford@housingmatters.org.uk") | Out-Null

# -- paragraph: "Housing Matters" --
$p1 = $d.Paragraphs(10)
$p1.Range.InsertParagraphAfter()
$p2 = $d.Paragraphs(11)
$d.Range($p2.Range.Start, $p2.Range.Start).InsertAfter("Housing Matters")

# -- paragraph: hyperlink to https://housingmatters.org.uk/ --
$p2 = $d.Paragraphs(11)
$p2.Range.InsertParagraphAfter()
$p3 = $d.Paragraphs(12)
$d.Range($p3.Range.Start, $p3.Range.Start).InsertAfter("X ")
$p3 = $d.Paragraphs(12)
$urlSel = $d.Range($p3.Range.Start, $p3.Range.Start + 1)
$d.Hyperlinks.Add($urlSel, "https://housingmatters.org.uk/", $null, $null, "https://housingmatters.org.uk/") | Out-Null

# ---------------------------------------------------------------------
# 3) Materialise the Hyperlink / Unresolved Mention character styles
# ---------------------------------------------------------------------
$hlStyle = $d.Styles.Add("Hyperlink", 2)
$hlStyle.BaseStyle = "DefaultParagraphFont"
$hlStyle.Priority = 99
$hlStyle.UnhideWhenUsed = $true
$hlStyle.Font.Underline = 1
$hlStyle.Font.TextColor.ObjectThemeColor = 10

$umStyle = $d.Styles.Add("Unresolved Mention", 2)
$umStyle.BaseStyle = "DefaultParagraphFont"
$umStyle.Priority = 99
$umStyle.UnhideWhenUsed = $true
$umStyle.Font.Color = 6053472

Write-Output "done"
